# Apply numeric corrections to Leve profit calculations across sheets.
# Each block targets one worksheet by name and updates specific cells
# (Cells.Item(row, col)) to match the refreshed market-price snapshot.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(15, 8).Value = 39.3
$ws.Cells.Item(15, 9).Value = 39.3
$ws.Cells.Item(15, 11).Value = 117.9
$ws.Cells.Item(15, 13).Value = 51.10000000000001
$ws.Cells.Item(64, 8).Value = 65624.81
$ws.Cells.Item(64, 9).Value = 252274.5
$ws.Cells.Item(64, 10).Value = 3408.25
$ws.Cells.Item(64, 11).Value = 252274.5
$ws.Cells.Item(64, 12).Value = 3408.25
$ws.Cells.Item(64, 13).Value = -252026.5
$ws.Cells.Item(64, 14).Value = -3904.25
$ws.Cells.Item(67, 8).Value = 65624.81
$ws.Cells.Item(67, 9).Value = 252274.5
$ws.Cells.Item(67, 10).Value = 3408.25
$ws.Cells.Item(67, 11).Value = 252274.5
$ws.Cells.Item(67, 12).Value = 3408.25
$ws.Cells.Item(67, 13).Value = -251416.5
$ws.Cells.Item(67, 14).Value = -5124.25
$ws.Cells.Item(70, 8).Value = 1110
$ws.Cells.Item(70, 9).Value = 998.5714
$ws.Cells.Item(70, 10).Value = 1500
$ws.Cells.Item(70, 11).Value = 2995.7142
$ws.Cells.Item(70, 12).Value = 4500
$ws.Cells.Item(70, 13).Value = -2725.7142
$ws.Cells.Item(70, 14).Value = -5040
$ws.Cells.Item(73, 8).Value = 1110
$ws.Cells.Item(73, 9).Value = 998.5714
$ws.Cells.Item(73, 10).Value = 1500
$ws.Cells.Item(73, 11).Value = 2995.7142
$ws.Cells.Item(73, 12).Value = 4500
$ws.Cells.Item(73, 13).Value = -2059.7142
$ws.Cells.Item(73, 14).Value = -6372
$ws.Cells.Item(76, 8).Value = 2996.5356
$ws.Cells.Item(76, 9).Value = 2980.8845
$ws.Cells.Item(76, 10).Value = 3200
$ws.Cells.Item(76, 11).Value = 2980.8845
$ws.Cells.Item(76, 12).Value = 3200
$ws.Cells.Item(76, 13).Value = -2665.8845
$ws.Cells.Item(76, 14).Value = -3830
$ws.Cells.Item(79, 8).Value = 2996.5356
$ws.Cells.Item(79, 9).Value = 2980.8845
$ws.Cells.Item(79, 10).Value = 3200
$ws.Cells.Item(79, 11).Value = 2980.8845
$ws.Cells.Item(79, 12).Value = 3200
$ws.Cells.Item(79, 13).Value = -1888.8845
$ws.Cells.Item(79, 14).Value = -5384
$ws.Cells.Item(82, 8).Value = 200004600
$ws.Cells.Item(82, 9).Value = 5760.5
$ws.Cells.Item(82, 10).Value = 1000000000
$ws.Cells.Item(82, 11).Value = 17281.5
$ws.Cells.Item(82, 12).Value = 3000000000
$ws.Cells.Item(82, 13).Value = -16875.5
$ws.Cells.Item(82, 14).Value = -3000000812
$ws.Cells.Item(85, 8).Value = 200004600
$ws.Cells.Item(85, 9).Value = 5760.5
$ws.Cells.Item(85, 10).Value = 1000000000
$ws.Cells.Item(85, 11).Value = 17281.5
$ws.Cells.Item(85, 12).Value = 3000000000
$ws.Cells.Item(85, 13).Value = -15877.5
$ws.Cells.Item(85, 14).Value = -3000002808
$ws.Cells.Item(88, 8).Value = 22675.334
$ws.Cells.Item(88, 9).Value = 8210
$ws.Cells.Item(88, 10).Value = 95002
$ws.Cells.Item(88, 11).Value = 8210
$ws.Cells.Item(88, 12).Value = 95002
$ws.Cells.Item(88, 13).Value = -7804
$ws.Cells.Item(88, 14).Value = -95814
$ws.Cells.Item(91, 8).Value = 22675.334
$ws.Cells.Item(91, 9).Value = 8210
$ws.Cells.Item(91, 10).Value = 95002
$ws.Cells.Item(91, 11).Value = 8210
$ws.Cells.Item(91, 12).Value = 95002
$ws.Cells.Item(91, 13).Value = -6806
$ws.Cells.Item(91, 14).Value = -97810
$ws.Cells.Item(97, 8).Value = 85608.38
$ws.Cells.Item(97, 9).Value = 633.3333
$ws.Cells.Item(97, 10).Value = 111100.9
$ws.Cells.Item(97, 11).Value = 1899.9999
$ws.Cells.Item(97, 12).Value = 333302.7
$ws.Cells.Item(97, 13).Value = -1403.9999
$ws.Cells.Item(97, 14).Value = -334294.7
$ws.Cells.Item(100, 8).Value = 4072.7273
$ws.Cells.Item(100, 9).Value = 4760
$ws.Cells.Item(100, 10).Value = 3500
$ws.Cells.Item(100, 11).Value = 4760
$ws.Cells.Item(100, 12).Value = 3500
$ws.Cells.Item(100, 13).Value = -4219
$ws.Cells.Item(100, 14).Value = -4582
$ws.Cells.Item(103, 8).Value = 349
$ws.Cells.Item(103, 9).Value = 349
$ws.Cells.Item(103, 10).Value = 0
$ws.Cells.Item(103, 11).Value = 1047
$ws.Cells.Item(103, 12).Value = 0
$ws.Cells.Item(103, 13).Value = -461
$ws.Cells.Item(103, 14).ClearContents()
$ws.Cells.Item(106, 8).Value = 79535.7
$ws.Cells.Item(106, 9).Value = 1953.3334
$ws.Cells.Item(106, 11).Value = 1953.3334
$ws.Cells.Item(106, 13).Value = -1322.3334
$ws.Cells.Item(112, 8).Value = 1073.8
$ws.Cells.Item(112, 10).Value = 1085.1538
$ws.Cells.Item(112, 12).Value = 3255.4614
$ws.Cells.Item(112, 14).Value = -5471.4614
$ws.Cells.Item(115, 8).Value = 1917.091
$ws.Cells.Item(115, 9).Value = 1917.091
$ws.Cells.Item(115, 11).Value = 5751.272999999999
$ws.Cells.Item(115, 13).Value = -4184.272999999999
$ws.Cells.Item(118, 8).Value = 1950
$ws.Cells.Item(118, 10).Value = 0
$ws.Cells.Item(118, 12).Value = 0
$ws.Cells.Item(118, 14).ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 2023.234
$ws.Cells.Item(61, 9).Value = 1861.5358
$ws.Cells.Item(61, 10).Value = 2261.5264
$ws.Cells.Item(61, 11).Value = 1861.5358
$ws.Cells.Item(61, 12).Value = 2261.5264
$ws.Cells.Item(61, 13).Value = -1649.5358
$ws.Cells.Item(61, 14).Value = -2685.5264
$ws.Cells.Item(132, 8).Value = 16669164
$ws.Cells.Item(132, 9).Value = 25001982
$ws.Cells.Item(132, 10).Value = 3529.4
$ws.Cells.Item(132, 11).Value = 75005946
$ws.Cells.Item(132, 12).Value = 10588.2
$ws.Cells.Item(132, 13).Value = -75003416
$ws.Cells.Item(132, 14).Value = -15648.2
$ws.Cells.Item(136, 8).Value = 2023.234
$ws.Cells.Item(136, 9).Value = 1861.5358
$ws.Cells.Item(136, 10).Value = 2261.5264
$ws.Cells.Item(136, 11).Value = 5584.607400000001
$ws.Cells.Item(136, 12).Value = 6784.5792
$ws.Cells.Item(136, 13).Value = -3034.607400000001
$ws.Cells.Item(136, 14).Value = -11884.5792

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(68, 8).Value = 68571.42999999999
$ws.Cells.Item(68, 10).Value = 68571.42999999999
$ws.Cells.Item(68, 12).Value = 68571.42999999999
$ws.Cells.Item(68, 14).Value = -70069.42999999999
$ws.Cells.Item(71, 8).Value = 68571.42999999999
$ws.Cells.Item(71, 10).Value = 68571.42999999999
$ws.Cells.Item(71, 12).Value = 205714.29
$ws.Cells.Item(71, 14).Value = -213202.29

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(130, 8).Value = 1538.125
$ws.Cells.Item(130, 10).Value = 1670.7142
$ws.Cells.Item(130, 12).Value = 5012.142599999999
$ws.Cells.Item(130, 14).Value = -15052.1426
$ws.Cells.Item(131, 8).Value = 32562.91
$ws.Cells.Item(131, 9).Value = 99999
$ws.Cells.Item(131, 10).Value = 31541.152
$ws.Cells.Item(131, 11).Value = 299997
$ws.Cells.Item(131, 12).Value = 94623.45599999999
$ws.Cells.Item(131, 13).Value = -294957
$ws.Cells.Item(131, 14).Value = -104703.456

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(132, 8).Value = 3613.6511
$ws.Cells.Item(132, 9).Value = 3485.1765
$ws.Cells.Item(132, 10).Value = 4099
$ws.Cells.Item(132, 11).Value = 10455.5295
$ws.Cells.Item(132, 12).Value = 12297
$ws.Cells.Item(132, 13).Value = -7925.529500000001
$ws.Cells.Item(132, 14).Value = -17357

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(100, 8).Value = 1818
$ws.Cells.Item(100, 9).Value = 1616.6471
$ws.Cells.Item(100, 10).Value = 2673.75
$ws.Cells.Item(100, 11).Value = 1616.6471
$ws.Cells.Item(100, 12).Value = 2673.75
$ws.Cells.Item(100, 13).Value = -1075.6471
$ws.Cells.Item(100, 14).Value = -3755.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132, 8).Value = 1518.9459
$ws.Cells.Item(132, 9).Value = 1103.4333
$ws.Cells.Item(132, 10).Value = 3299.7144
$ws.Cells.Item(132, 11).Value = 3310.2999
$ws.Cells.Item(132, 12).Value = 9899.143199999999
$ws.Cells.Item(132, 13).Value = -780.2999
$ws.Cells.Item(132, 14).Value = -14959.1432
